$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:56:54"
$ws.Cells.Item(3, 1).Value = "Total filas: 243"
$ws.Cells.Item(54, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(55, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(56, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(65, 1).Value = "07:18:07"
$ws.Cells.Item(65, 3).Value = "215D_EL PATO"
$ws.Cells.Item(65, 4).Value = 34
$ws.Cells.Item(66, 1).Value = "07:50:27"
$ws.Cells.Item(66, 3).Value = "10_OLMOS"
$ws.Cells.Item(66, 4).Value = 2
$ws.Cells.Item(89, 3).Value = "10_OLMOS"
$ws.Cells.Item(90, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(99, 1).Value = "07:18:07"
$ws.Cells.Item(99, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(99, 4).Value = 113
$ws.Cells.Item(100, 1).Value = "07:50:27"
$ws.Cells.Item(100, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(100, 4).Value = 81
$ws.Cells.Item(108, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(110, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(148, 3).Value = "215C_EL PATO"
$ws.Cells.Item(149, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(150, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(190, 3).Value = "215A_EL PATO"
$ws.Cells.Item(191, 3).Value = "14_ABASTO"
$ws.Cells.Item(208, 1).Value = "12:56:54"
$ws.Cells.Item(208, 2).Value = "12:56"
$ws.Cells.Item(208, 3).Value = "10_OLMOS"
$ws.Cells.Item(208, 4).Value = 0
$ws.Cells.Item(209, 1).Value = "12:56:54"
$ws.Cells.Item(209, 2).Value = "13:02"
$ws.Cells.Item(209, 3).Value = "15_ABASTO"
$ws.Cells.Item(209, 4).Value = 6
$ws.Cells.Item(210, 2).Value = "13:03"
$ws.Cells.Item(210, 3).Value = "14_ABASTO"
$ws.Cells.Item(210, 4).Value = 23
$ws.Cells.Item(211, 1).Value = "12:40:02"
$ws.Cells.Item(211, 2).Value = "13:05"
$ws.Cells.Item(211, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(211, 4).Value = 25
$ws.Cells.Item(212, 1).Value = "12:56:54"
$ws.Cells.Item(212, 2).Value = "13:06"
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(213, 1).Value = "12:56:54"
$ws.Cells.Item(213, 2).Value = "13:07"
$ws.Cells.Item(213, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(213, 4).Value = 11
$ws.Cells.Item(214, 1).Value = "12:56:54"
$ws.Cells.Item(214, 2).Value = "13:07"
$ws.Cells.Item(214, 3).Value = "10_OLMOS"
$ws.Cells.Item(214, 4).Value = 11
$ws.Cells.Item(215, 1).Value = "12:40:02"
$ws.Cells.Item(215, 2).Value = "13:07"
$ws.Cells.Item(215, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(215, 4).Value = 27
$ws.Cells.Item(216, 2).Value = "13:08"
$ws.Cells.Item(216, 3).Value = "10_OLMOS"
$ws.Cells.Item(216, 4).Value = 28
$ws.Cells.Item(217, 1).Value = "12:56:54"
$ws.Cells.Item(217, 2).Value = "13:08"
$ws.Cells.Item(217, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(217, 4).Value = 12
$ws.Cells.Item(218, 1).Value = "12:40:02"
$ws.Cells.Item(218, 2).Value = "13:09"
$ws.Cells.Item(218, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(218, 4).Value = 29
$ws.Cells.Item(219, 1).Value = "12:56:54"
$ws.Cells.Item(219, 2).Value = "13:13"
$ws.Cells.Item(219, 3).Value = "215D_EL PATO"
$ws.Cells.Item(219, 4).Value = 17
$ws.Cells.Item(220, 2).Value = "13:14"
$ws.Cells.Item(220, 3).Value = "215D_EL PATO"
$ws.Cells.Item(220, 4).Value = 34
$ws.Cells.Item(221, 1).Value = "12:56:54"
$ws.Cells.Item(221, 2).Value = "13:14"
$ws.Cells.Item(221, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(221, 4).Value = 18
$ws.Cells.Item(222, 1).Value = "12:40:02"
$ws.Cells.Item(222, 2).Value = "13:15"
$ws.Cells.Item(222, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(222, 4).Value = 35
$ws.Cells.Item(223, 1).Value = "12:56:54"
$ws.Cells.Item(223, 2).Value = "13:19"
$ws.Cells.Item(223, 3).Value = "10_OLMOS"
$ws.Cells.Item(223, 4).Value = 23
$ws.Cells.Item(224, 1).Value = "12:56:54"
$ws.Cells.Item(224, 2).Value = "13:20"
$ws.Cells.Item(224, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(224, 4).Value = 24
$ws.Cells.Item(225, 1).Value = "12:40:02"
$ws.Cells.Item(225, 2).Value = "13:20"
$ws.Cells.Item(225, 3).Value = "10_OLMOS"
$ws.Cells.Item(225, 4).Value = 40
$ws.Cells.Item(226, 2).Value = "13:21"
$ws.Cells.Item(226, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(226, 4).Value = 41
$ws.Cells.Item(227, 1).Value = "12:56:54"
$ws.Cells.Item(227, 2).Value = "13:26"
$ws.Cells.Item(227, 3).Value = "14_ABASTO"
$ws.Cells.Item(227, 4).Value = 30
$ws.Cells.Item(228, 1).Value = "12:56:54"
$ws.Cells.Item(228, 2).Value = "13:26"
$ws.Cells.Item(228, 3).Value = "15_ABASTO"
$ws.Cells.Item(228, 4).Value = 30
$ws.Cells.Item(229, 2).Value = "13:27"
$ws.Cells.Item(229, 3).Value = "14_ABASTO"
$ws.Cells.Item(229, 4).Value = 47
$ws.Cells.Item(230, 1).Value = "12:56:54"
$ws.Cells.Item(230, 2).Value = "13:33"
$ws.Cells.Item(230, 3).Value = "10_OLMOS"
$ws.Cells.Item(230, 4).Value = 37
$ws.Cells.Item(231, 1).Value = "12:56:54"
$ws.Cells.Item(231, 2).Value = "13:34"
$ws.Cells.Item(231, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(231, 4).Value = 38
$ws.Cells.Item(232, 1).Value = "12:56:54"
$ws.Cells.Item(232, 2).Value = "13:46"
$ws.Cells.Item(232, 3).Value = "17_ROMERO"
$ws.Cells.Item(232, 4).Value = 50
$ws.Cells.Item(233, 1).Value = "12:56:54"
$ws.Cells.Item(233, 2).Value = "13:46"
$ws.Cells.Item(233, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(233, 4).Value = 50
$ws.Cells.Item(234, 1).Value = "12:56:54"
$ws.Cells.Item(234, 2).Value = "13:50"
$ws.Cells.Item(234, 3).Value = "215A_EL PATO"
$ws.Cells.Item(234, 4).Value = 54
$ws.Cells.Item(234, 5).Value = "LP1912"
$ws.Cells.Item(235, 1).Value = "12:56:54"
$ws.Cells.Item(235, 2).Value = "13:50"
$ws.Cells.Item(235, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(235, 4).Value = 54
$ws.Cells.Item(235, 5).Value = "LP1912"
$ws.Cells.Item(236, 1).Value = "12:40:02"
$ws.Cells.Item(236, 2).Value = "13:51"
$ws.Cells.Item(236, 3).Value = "215A_EL PATO"
$ws.Cells.Item(236, 4).Value = 71
$ws.Cells.Item(236, 5).Value = "LP1912"
$ws.Cells.Item(237, 1).Value = "12:56:54"
$ws.Cells.Item(237, 2).Value = "13:55"
$ws.Cells.Item(237, 3).Value = "225_GOMEZ"
$ws.Cells.Item(237, 4).Value = 59
$ws.Cells.Item(237, 5).Value = "LP1912"
$ws.Cells.Item(238, 1).Value = "12:40:02"
$ws.Cells.Item(238, 2).Value = "13:56"
$ws.Cells.Item(238, 3).Value = "225_GOMEZ"
$ws.Cells.Item(238, 4).Value = 76
$ws.Cells.Item(238, 5).Value = "LP1912"
$ws.Cells.Item(239, 1).Value = "12:56:54"
$ws.Cells.Item(239, 2).Value = "13:56"
$ws.Cells.Item(239, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(239, 4).Value = 60
$ws.Cells.Item(239, 5).Value = "LP1912"
$ws.Cells.Item(240, 1).Value = "12:40:02"
$ws.Cells.Item(240, 2).Value = "13:57"
$ws.Cells.Item(240, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(240, 4).Value = 77
$ws.Cells.Item(240, 5).Value = "LP1912"
$ws.Cells.Item(241, 1).Value = "12:56:54"
$ws.Cells.Item(241, 2).Value = "14:04"
$ws.Cells.Item(241, 3).Value = "17_ROMERO"
$ws.Cells.Item(241, 4).Value = 68
$ws.Cells.Item(241, 5).Value = "LP1912"
$ws.Cells.Item(242, 1).Value = "12:56:54"
$ws.Cells.Item(242, 2).Value = "14:16"
$ws.Cells.Item(242, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(242, 4).Value = 80
$ws.Cells.Item(242, 5).Value = "LP1912"
$ws.Cells.Item(243, 1).Value = "12:40:02"
$ws.Cells.Item(243, 2).Value = "14:17"
$ws.Cells.Item(243, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(243, 4).Value = 97
$ws.Cells.Item(243, 5).Value = "LP1912"
$ws.Cells.Item(244, 1).Value = "12:56:54"
$ws.Cells.Item(244, 2).Value = "14:19"
$ws.Cells.Item(244, 3).Value = "215C_EL PATO"
$ws.Cells.Item(244, 4).Value = 83
$ws.Cells.Item(244, 5).Value = "LP1912"
$ws.Cells.Item(245, 1).Value = "12:56:54"
$ws.Cells.Item(245, 2).Value = "14:20"
$ws.Cells.Item(245, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(245, 4).Value = 84
$ws.Cells.Item(245, 5).Value = "LP1912"
$ws.Cells.Item(246, 1).Value = "12:40:02"
$ws.Cells.Item(246, 2).Value = "14:20"
$ws.Cells.Item(246, 3).Value = "215C_EL PATO"
$ws.Cells.Item(246, 4).Value = 100
$ws.Cells.Item(246, 5).Value = "LP1912"
$ws.Cells.Item(247, 1).Value = "12:40:02"
$ws.Cells.Item(247, 2).Value = "14:21"
$ws.Cells.Item(247, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(247, 4).Value = 101
$ws.Cells.Item(247, 5).Value = "LP1912"
$ws.Cells.Item(248, 1).Value = "12:56:54"
$ws.Cells.Item(248, 2).Value = "14:49"
$ws.Cells.Item(248, 3).Value = "14_ABASTO"
$ws.Cells.Item(248, 4).Value = 113
$ws.Cells.Item(248, 5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:56:54"
$ws.Cells.Item(3, 1).Value = "Total filas: 26"
$ws.Cells.Item(26, 1).Value = "12:56:54"
$ws.Cells.Item(26, 4).Value = 17
$ws.Cells.Item(28, 1).Value = "12:56:54"
$ws.Cells.Item(28, 4).Value = 54
$ws.Cells.Item(30, 1).Value = "12:56:54"
$ws.Cells.Item(30, 2).Value = "14:19"
$ws.Cells.Item(30, 4).Value = 83
$ws.Cells.Item(31, 1).Value = "12:40:02"
$ws.Cells.Item(31, 2).Value = "14:20"
$ws.Cells.Item(31, 3).Value = "215C_EL PATO"
$ws.Cells.Item(31, 4).Value = 100
$ws.Cells.Item(31, 5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:56:54"
$ws.Cells.Item(3, 1).Value = "Total filas: 30"
$ws.Cells.Item(32, 1).Value = "12:56:54"
$ws.Cells.Item(32, 4).Value = 34
$ws.Cells.Item(34, 1).Value = "12:56:54"
$ws.Cells.Item(34, 4).Value = 73
$ws.Cells.Item(35, 1).Value = "12:56:54"
$ws.Cells.Item(35, 2).Value = "14:52"
$ws.Cells.Item(35, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(35, 4).Value = 116
$ws.Cells.Item(35, 5).Value = "L6203"
